{"js": "// Remove the \"Docente(s) Respons\u00e1vel(eis)\" bullet-list paragraph that names\n// the responsible professor (1097178 - Jo\u00e3o Batista de Almeida e Silva).\n// The heading paragraph itself (\"Docente(s) Respons\u00e1vel(eis) \") is kept;\n// only the following List Bullet paragraph is deleted entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"1097178 - Jo\u00e3o Batista de Almeida e Silva\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === targetText) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Docente(s) Respons\u00e1vel(eis)\" bullet-list paragraph that names\n# the responsible professor (1097178 - Jo\u00e3o Batista de Almeida e Silva).\n# The heading paragraph itself (\"Docente(s) Respons\u00e1vel(eis) \") is kept;\n# only the following List Bullet paragraph is deleted entirely.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"1097178 - Jo\u00e3o Batista de Almeida e Silva\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
